$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.965.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.076.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.60%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.67%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.068.41"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.82%  "
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("E10").Value = "  +5.27%  "
$ws.Range("E11").Value = "  +11.91%  "
$ws.Range("E12").Value = "  +2.13%  "
$ws.Range("E13").Value = "  +4.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.37%  "
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.588.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.077.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "61.878.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "448.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.06%  "
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("E23").Value = "  +4.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("E27").Value = "  +4.19%  "
$ws.Range("B28").Value = "FirstDigitalUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.26%  "
$ws.Range("E30").Value = "  +4.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E32").Value = "  +12.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.41%  "
$ws.Range("E34").Value = "  +5.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0793"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.33%  "
$ws.Range("E36").Value = "  +2.67%  "
$ws.Range("E37").Value = "  +4.88%  "
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "420.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.45%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.928.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.67%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0371"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.49%  "
$ws.Range("E44").Value = "  +9.41%  "
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("E46").Value = "  +6.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.14%  "
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("E51").Value = "  +4.01%  "
